$d = $word.ActiveDocument

# The parser now uses TokenIteratorFieldRewriterSplit: the two "malformed"
# fields in this doc (an `m:for` missing its `|` and the matching
# `m:endfor`) are no longer real Word fields (w:fldChar/w:instrText) -
# they become plain literal text runs, brace-delimited, so the token
# iterator re-splits them itself:
#   { m:for v     self.eAllStructuralFeatures }   ->  {m:|for v| |self.eAllStructuralFeatures}
#   { m:endfor }                                   ->  {|m:|endfor}
#
# Find the paragraph holding those fields (still real fields in the
# starting document) and rebuild it as plain text runs.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find the paragraph containing the malformed fields."
}

# Sanity-check we found the expected two fields before nuking them.
$fields = $targetPara.Range.Fields
if ($fields.Count -ne 2) {
    throw "Expected 2 fields in the target paragraph, found $($fields.Count)."
}

# Runs that stay exactly as-is (everything between/after the two fields):
#   "    "  (4 spaces) / the red "<---Malformed tag m:for, no '|' found."
#   / "A paragraph" / the _GoBack bookmark.
# Runs that used to be w:fldChar/w:instrText become plain w:t runs:
#   field 1 "m:for v     self.eAllStructuralFeatures " -> {m:  /  for v  /  (sp)  /  self.eAllStructuralFeatures}
#   field 2 " m:endfor "                                -> {    /  m:      /  endfor}
$run_field1_part1 = '<w:r><w:t>{m:</w:t></w:r>'
$run_field1_part2 = '<w:r><w:t>for v</w:t></w:r>'
$run_field1_part3 = '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$run_field1_part4 = '<w:r><w:t xml:space="preserve">self.eAllStructuralFeatures}</w:t></w:r>'

$run_spaces      = '<w:r><w:t xml:space="preserve">    </w:t></w:r>'
$run_arrow       = '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>'
$run_malformed   = '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>Malformed tag m:for, no ' + "'" + '|' + "'" + ' found.</w:t></w:r>'
$run_aparagraph  = '<w:r w:rsidR="00146761"><w:t>A paragraph</w:t></w:r>'
$bookmark        = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$run_field2_part1 = '<w:r><w:t>{</w:t></w:r>'
$run_field2_part2 = '<w:r><w:t>m:</w:t></w:r>'
$run_field2_part3 = '<w:r><w:t xml:space="preserve">endfor}</w:t></w:r>'

$newParagraphInner = $run_field1_part1 + $run_field1_part2 + $run_field1_part3 + $run_field1_part4 + `
                      $run_spaces + $run_arrow + $run_malformed + $run_aparagraph + $bookmark + `
                      $run_field2_part1 + $run_field2_part2 + $run_field2_part3

$newParagraphXml = '<w:p w:rsidR="007A2DC4" w:rsidRDefault="00B31BB7">' + $newParagraphInner + '</w:p>'

$packageXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + $newParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replace the whole paragraph (this engine treats any range touching a
# w:fldChar as "replace the owning paragraph", so we target the full
# paragraph range and reconstruct it rather than surgically poking at the
# field runs in place).
$targetPara.Range.InsertXML($packageXml)
